$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the computed "Utility" values in column D (re-run of the SeeDB
# utility scoring for this view, early-August numbers).
$ws.Cells.Item(2, 4).Value = 0.07504891401632832193
$ws.Cells.Item(3, 4).Value = 0.06819214133962332725
$ws.Cells.Item(4, 4).Value = 0.05671152682990233257
$ws.Cells.Item(5, 4).Value = 0.05191711710604315216
$ws.Cells.Item(6, 4).Value = 0.05001535590561854983
$ws.Cells.Item(7, 4).Value = 0.04815802989787729138
$ws.Cells.Item(11, 4).Value = 0.45719128658069091431
$ws.Cells.Item(12, 4).Value = 0.17580437891423281171
$ws.Cells.Item(13, 4).Value = 0.14463815981609828842
$ws.Cells.Item(14, 4).Value = 0.13320909625701268797
$ws.Cells.Item(15, 4).Value = 0.06357015097151620664
$ws.Cells.Item(16, 4).Value = 0.04740627170963671727
$ws.Cells.Item(17, 4).Value = 0.02457913499594240048
$ws.Cells.Item(18, 4).Value = 0.01874892616203872875
$ws.Cells.Item(19, 4).Value = 0.01089977047417665067
$ws.Cells.Item(20, 4).Value = 0.05953528455220139654
$ws.Cells.Item(21, 4).Value = 0.05764001202723850747
$ws.Cells.Item(23, 4).Value = 0.03443125109773714776
$ws.Cells.Item(24, 4).Value = 0.02480245090811101019
$ws.Cells.Item(25, 4).Value = 0.02233237129794678894
$ws.Cells.Item(26, 4).Value = 0.02034543796873626922

# Reset the view: scroll back to the top-left and move the selection.
$window = $excel.ActiveWindow
$window.ScrollRow = 1
$window.ScrollColumn = 1
$ws.Range("K16").Select()
